$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Registros")

# B24 holds a numeric-looking score but must remain text (as in the rest of
# column B), so force text formatting while writing it, then restore the
# original (General) number format.
$cellB24 = $ws.Range("B24")
$origFormat = $cellB24.NumberFormat
$cellB24.NumberFormat = "@"
$cellB24.Value = "5"
$cellB24.NumberFormat = $origFormat

$ws.Range("C24").Value = "none"
$ws.Range("D24").Value = "b"
$ws.Range("E24").Value = "c"
$ws.Range("F24").Value = "d"
$ws.Range("G24").Value = "e"
$ws.Range("H24").Value = "d"
$ws.Range("I24").Value = "c"
$ws.Range("J24").Value = "b"
$ws.Range("K24").Value = "a"
$ws.Range("L24").Value = "b"
$ws.Range("M24").Value = "c"
$ws.Range("N24").Value = "d"
$ws.Range("O24").Value = "e"
$ws.Range("P24").Value = "d"
$ws.Range("Q24").Value = "c"
$ws.Range("R24").Value = "b"
$ws.Range("S24").Value = "a"
$ws.Range("T24").Value = "b"
$ws.Range("U24").Value = "c"
$ws.Range("V24").Value = "d"
$ws.Range("W24").Value = "e"
